$wb = $excel.ActiveWorkbook

# --- "Object Code" sheet: add particle-effect / ground-mask object codes ---
$objectCode = $wb.Worksheets.Item("Object Code")
$objectCode.Activate()
$objectCode.Range("B5").Value = "particle effect"
$objectCode.Range("B6").Value = "ground mask"
$objectCode.Range("B6").Select() | Out-Null

# --- "Image ID" sheet: add ammo / health / boom image ids, becomes the active sheet ---
$imageId = $wb.Worksheets.Item("Image ID")
$imageId.Activate()
$imageId.Range("A7").Value = 5
$imageId.Range("B7").Value = "ammo"
$imageId.Range("A8").Value = 6
$imageId.Range("B8").Value = "health"
$imageId.Range("A9").Value = 7
$imageId.Range("B9").Value = "boom"
$imageId.Range("B9").Select() | Out-Null
